$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("BD3").Value = 151

# Row 4
$ws.Range("G4").Value  = 2.8
$ws.Range("I4").Value  = 2.5
$ws.Range("J4").Value  = 3.4
$ws.Range("L4").Value  = 3.1
$ws.Range("Y4").Value  = 11
$ws.Range("AA4").Value = 23
$ws.Range("AI4").Value = 12
$ws.Range("AK4").Value = 23
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 21
$ws.Range("BD4").Value = 151

# Row 5
$ws.Range("G5").Value  = 2.7
$ws.Range("I5").Value  = 2.55
$ws.Range("J5").Value  = 3.25
$ws.Range("L5").Value  = 3.1
$ws.Range("W5").Value  = 10
$ws.Range("X5").Value  = 15
$ws.Range("AL5").Value = 19
$ws.Range("AM5").Value = 26
$ws.Range("AQ5").Value = 51
$ws.Range("AX5").Value = 13
$ws.Range("AY5").Value = 21
